{"js": "// The Bibliografia paragraph holds a single run whose text is all five\n// numbered references run together. Split it into five references\n// separated by manual line breaks (Word's \"\\u000B\" line-break\n// character), so OOXML serializes the run as\n// <w:r><w:t>ref1</w:t><w:br/><w:t>ref2</w:t><w:br/>...</w:r> instead\n// of one giant <w:t>.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the bibliography paragraph by its distinctive leading text\n// (robust against any unrelated structural changes elsewhere in the\n// document, rather than relying on a hard-coded paragraph index).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text || \"\";\n  if (t.indexOf(\"[1] Peddy\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Bibliography paragraph not found\");\n}\n\n// The five individual reference strings, split out of the original\n// concatenated text (kept verbatim, only the split points differ).\nconst refs = [\n  \"[1] Peddy, S. The art of mentoring \\u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\",\n  \"[2] Zachary, L. J. The Mentor\\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\\u00e7\\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\\u00e7\\u00e3o e Tutorado: oito anos a promover a integra\\u00e7\\u00e3o e o sucesso acad\\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\",\n  \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\",\n  \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\",\n  \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\\u00e7\\u00e3o em Engenharia. Minist\\u00e9rio da Educa\\u00e7\\u00e3o. CNE/CES, 2019.\",\n];\n\n// Empty the paragraph, then retype it as one run with embedded manual\n// line breaks between references.\ntarget.clear();\nawait context.sync();\n\ntarget.insertText(refs.join(\"\\u000B\"), \"Start\");\nawait context.sync();\n", "ps1": "# Split the single run in the Bibliografia paragraph into a run that\n# carries five reference strings separated by manual line breaks\n# (Word's vertical-tab line-break character, `v / Chr(11)), one break\n# between each of the five numbered references ([1]..[5]). OOXML\n# serializes each manual line break as <w:br/>, so the paragraph ends\n# up as <w:r><w:t>ref1</w:t><w:br/><w:t>ref2</w:t><w:br/>...</w:r>.\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography paragraph via Find (robust against any\n# unrelated structural changes elsewhere in the document), then expand\n# the found range to the whole paragraph so we can overwrite its text\n# in one shot.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"[1] Peddy\")\nif (-not $found) {\n    throw \"Bibliography paragraph not found\"\n}\n$rng.Expand(4)  # wdParagraph\n\n# The five individual reference strings, split out of the original\n# concatenated text (kept verbatim, only the split points differ).\n$ref1 = \"[1] Peddy, S. The art of mentoring \" + [char]0x2013 + \" Lead, follow and get out of the way. Houston: Bullion Books, 2001.\"\n$ref2 = \"[2] Zachary, L. J. The Mentor\" + [char]0x2019 + \"s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\" + [char]0x00E7 + [char]0x00E3 + \"o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\" + [char]0x00E7 + [char]0x00E3 + \"o e Tutorado: oito anos a promover a integra\" + [char]0x00E7 + [char]0x00E3 + \"o e o sucesso acad\" + [char]0x00E9 + \"mico no IST. Lisboa: IST Press, 2011. p. 19-27.\"\n$ref3 = \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\"\n$ref4 = \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\"\n$ref5 = \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\" + [char]0x00E7 + [char]0x00E3 + \"o em Engenharia. Minist\" + [char]0x00E9 + \"rio da Educa\" + [char]0x00E7 + [char]0x00E3 + \"o. CNE/CES, 2019.\"\n\n$lineBreak = [char]11\n$rng.Text = $ref1 + $lineBreak + $ref2 + $lineBreak + $ref3 + $lineBreak + $ref4 + $lineBreak + $ref5\n"}
